$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Update time_taken (F column) timestamps on the "data" sheet ---
$ws.Range("F2").Value = "2021-10-05 14:21:54.544386"
$ws.Range("F3").Value = "2021-10-05 14:21:54.544396"
$ws.Range("F4").Value = "2021-10-05 14:21:54.544399"
$ws.Range("F5").Value = "2021-10-05 14:21:54.544402"
$ws.Range("F6").Value = "2021-10-05 14:21:54.544405"
$ws.Range("F7").Value = "2021-10-05 14:21:54.544408"
$ws.Range("F8").Value = "2021-10-05 14:21:54.544411"
$ws.Range("F9").Value = "2021-10-05 14:21:54.544413"
$ws.Range("F10").Value = "2021-10-05 14:21:54.544416"
$ws.Range("F11").Value = "2021-10-05 14:21:54.544419"
$ws.Range("F12").Value = "2021-10-05 14:21:54.544422"
$ws.Range("F13").Value = "2021-10-05 14:21:54.544424"
$ws.Range("F14").Value = "2021-10-05 14:21:54.544427"
$ws.Range("F15").Value = "2021-10-05 14:21:54.544430"
$ws.Range("F16").Value = "2021-10-05 14:21:54.544432"
$ws.Range("F17").Value = "2021-10-05 14:21:54.544435"
$ws.Range("F18").Value = "2021-10-05 14:21:54.544438"
$ws.Range("F19").Value = "2021-10-05 14:21:54.544441"
$ws.Range("F20").Value = "2021-10-05 14:21:54.544444"
$ws.Range("F21").Value = "2021-10-05 14:21:54.544447"
$ws.Range("F22").Value = "2021-10-05 14:21:54.544449"
$ws.Range("F23").Value = "2021-10-05 14:21:54.544452"
$ws.Range("F24").Value = "2021-10-05 14:21:54.544455"
$ws.Range("F25").Value = "2021-10-05 14:21:54.544458"
$ws.Range("F26").Value = "2021-10-05 14:21:54.544461"
$ws.Range("F27").Value = "2021-10-05 14:21:54.544464"
$ws.Range("F28").Value = "2021-10-05 14:21:54.544467"
$ws.Range("F29").Value = "2021-10-05 14:21:54.544469"
$ws.Range("F30").Value = "2021-10-05 14:21:54.544472"
$ws.Range("F31").Value = "2021-10-05 14:21:54.544475"
$ws.Range("F32").Value = "2021-10-05 14:21:54.544478"
$ws.Range("F33").Value = "2021-10-05 14:21:54.544481"
$ws.Range("F34").Value = "2021-10-05 14:21:54.544484"
$ws.Range("F35").Value = "2021-10-05 14:21:54.544487"
$ws.Range("F36").Value = "2021-10-05 14:21:54.544489"
$ws.Range("F37").Value = "2021-10-05 14:21:54.544492"
$ws.Range("F38").Value = "2021-10-05 14:21:54.544495"
$ws.Range("F39").Value = "2021-10-05 14:21:54.544497"
$ws.Range("F40").Value = "2021-10-05 14:21:54.544500"
$ws.Range("F41").Value = "2021-10-05 14:21:54.544503"
$ws.Range("F42").Value = "2021-10-05 14:21:54.544506"
$ws.Range("F43").Value = "2021-10-05 14:21:54.544509"
$ws.Range("F44").Value = "2021-10-05 14:21:54.544512"
$ws.Range("F45").Value = "2021-10-05 14:21:54.544515"
$ws.Range("F46").Value = "2021-10-05 14:21:54.544517"
$ws.Range("F47").Value = "2021-10-05 14:21:54.544520"
$ws.Range("F48").Value = "2021-10-05 14:21:54.544523"
$ws.Range("F49").Value = "2021-10-05 14:21:54.544525"
$ws.Range("F50").Value = "2021-10-05 14:21:54.544528"
$ws.Range("F51").Value = "2021-10-05 14:21:54.544531"
$ws.Range("F52").Value = "2021-10-05 14:21:54.544533"
$ws.Range("F53").Value = "2021-10-05 14:21:54.544536"
$ws.Range("F54").Value = "2021-10-05 14:21:54.544539"
$ws.Range("F55").Value = "2021-10-05 14:21:54.544542"
$ws.Range("F56").Value = "2021-10-05 14:21:54.544545"
$ws.Range("F57").Value = "2021-10-05 14:21:54.544548"
$ws.Range("F58").Value = "2021-10-05 14:21:54.544550"
$ws.Range("F59").Value = "2021-10-05 14:21:54.544553"
$ws.Range("F60").Value = "2021-10-05 14:21:54.544556"
$ws.Range("F61").Value = "2021-10-05 14:21:54.544559"
$ws.Range("F62").Value = "2021-10-05 14:21:54.544562"
$ws.Range("F63").Value = "2021-10-05 14:21:54.544564"

# --- Add the "metadata" worksheet after "data" ---
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Copy header styling (bold/border/centered) from the data sheet's header row,
# then overwrite with the metadata column headers.
$ws.Range("B1:F1").Copy($meta.Range("B1:F1"))
$ws.Range("F1").Copy($meta.Range("G1"))

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Copy the index-column style from data!A2 to metadata!A2, then set its value.
$ws.Range("A2").Copy($meta.Range("A2"))
$meta.Range("A2").Value = 0

$meta.Range("B2").Value = "Optic neuropathy"
$meta.Range("C2").Value = 186
$meta.Range("D2").Value = "'2.50"
$meta.Range("E2").Value = "2021-08-23T15:06:54.281695Z"
$meta.Range("F2").Value = "2021-10-05 14:21:54.540809"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/186/?format=json"

# Keep the "data" sheet as the active/selected tab (unchanged bookView).
$ws.Activate()
$ws.Range("A1").Select()
